# added 4wk low sales check
$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- "Forecast Comparison" sheet: recomputed Inventory Coverage (H),
#     Stockout Risk (I), Reorder Urgency (J) and Seasonality Index (L)
#     now that a 4-week low-sales check has been added. ---

# Row 2 (W10)
$wsForecast.Range("H2").Value = 7.5
$wsForecast.Range("L2").Value = 1.1

# Row 3 (W11)
$wsForecast.Range("H3").Value = 13
$wsForecast.Range("L3").Value = 1.05

# Row 4 (W12)
$wsForecast.Range("H4").Value = 6
$wsForecast.Range("L4").Value = 1.11

# Row 5 (W13)
$wsForecast.Range("H5").Value = 10
$wsForecast.Range("L5").Value = 1.02

# Row 6 (W14)
$wsForecast.Range("H6").Value = 9
$wsForecast.Range("J6").Value = "Normal"
$wsForecast.Range("L6").Value = 1.02

# Row 7 (W15)
$wsForecast.Range("H7").Value = 8
$wsForecast.Range("I7").Value = "Low"
$wsForecast.Range("J7").Value = "Normal"
$wsForecast.Range("L7").Value = 1.03

# Row 8 (W16)
$wsForecast.Range("H8").Value = 7
$wsForecast.Range("I8").Value = "Low"
$wsForecast.Range("J8").Value = "Normal"
$wsForecast.Range("L8").Value = 1

# Row 9 (W17)
$wsForecast.Range("H9").Value = 6
$wsForecast.Range("I9").Value = "Low"
$wsForecast.Range("J9").Value = "Normal"
$wsForecast.Range("L9").Value = 0.83

# Row 10 (W18)
$wsForecast.Range("H10").Value = 5
$wsForecast.Range("I10").Value = "Low"
$wsForecast.Range("J10").Value = "Normal"
$wsForecast.Range("L10").Value = 0.84

# Row 11 (W19)
$wsForecast.Range("H11").Value = 4
$wsForecast.Range("I11").Value = "Low"
$wsForecast.Range("J11").Value = "Normal"
$wsForecast.Range("L11").Value = 0.96

# Row 12 (W20)
$wsForecast.Range("H12").Value = 3
$wsForecast.Range("I12").Value = "Low"
$wsForecast.Range("J12").Value = "Normal"
$wsForecast.Range("L12").Value = 0.84

# Row 13 (W21)
$wsForecast.Range("H13").Value = 2
$wsForecast.Range("I13").Value = "Low"
$wsForecast.Range("J13").Value = "Normal"
$wsForecast.Range("L13").Value = 0.96

# Row 14 (W22)
$wsForecast.Range("H14").ClearContents()
$wsForecast.Range("I14").Value = "Low"
$wsForecast.Range("J14").Value = "Normal"
$wsForecast.Range("L14").Value = 0.87

# Row 15 (W23)
$wsForecast.Range("H15").Value = 1
$wsForecast.Range("I15").Value = "Low"
$wsForecast.Range("J15").Value = "Normal"

# Row 16 (W24)
$wsForecast.Range("L16").Value = 1.05

# Row 17 (W25)
$wsForecast.Range("H17").ClearContents()
$wsForecast.Range("I17").Value = "Low"
$wsForecast.Range("J17").Value = "Normal"
$wsForecast.Range("L17").Value = 1.03

# --- "Summary" sheet: forecast totals recalculated after the 4wk low
#     sales check was added. These are stored as text, like the rest
#     of the Value column, so force text with a leading apostrophe. ---

$wsSummary.Range("B9").Value  = "'3"   # Total Forecast (16 Weeks)
$wsSummary.Range("B10").Value = "'2"   # Total Forecast (8 Weeks)
$wsSummary.Range("B11").Value = "'1"   # Total Forecast (4 Weeks)
$wsSummary.Range("B12").Value = "'0"   # Max Forecast
$wsSummary.Range("B14").Value = "'0"   # Min Forecast
